# "cross in 3D view" - update TimeSheet hours/task text for Aban (Oct/Nov)
# and Bahman (Jan/Feb) sections, and rename a task description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- آبان 98 (Aban) section ---
# 3D View & Virtual Camera hours: 2 -> 3
$ws.Range("C63").Value = 3
# @Home hours for Aban: 2 -> 3 (keeps the SUM(C63:C63) total of 3 in sync)
$ws.Range("D66").Value = 3

# --- بهمن 98 (Bahman) section ---
# Rename task: "Get Tool/Ref Coordinates" -> "Get Tool/Ref Coordinates from Tracker"
$ws.Range("E84").Value = "* Get Tool/Ref Coordinates from Tracker"

# Hours updates
$ws.Range("C84").Value = 5
$ws.Range("C85").Value = 5
$ws.Range("C87").Value = 7
$ws.Range("C88").Value = 4
$ws.Range("C89").Value = 2

# @Parsiss / @Home totals for Bahman
$ws.Range("D91").Value = 5
$ws.Range("D92").Value = 22

# --- View / selection state ---
$ws.Range("E93").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 73
$win.ScrollColumn = 1
